# Refresh the crypto price-tracking table with newly scraped values.
# Source: "Updated symbol list on Thu Dec 22 10:11:37 UTC 2022 with GitHub Actions"
#
# The sheet's rank/volume labels in column E embed the coin's row-rank
# (e.g. "6KuCoinTokenKCS" for rank 6), so when a coin's rank shifts, its
# name (B), link (C), price (D) and label (E) all move together; the
# "Hora" column G is bumped from 9 to 10 for every row, and a handful of
# prices (D) are refreshed even where the rank didn't change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry is a single-cell text update. All of these columns hold
# plain text in the workbook (prices included, e.g. "0.001590" with a
# significant trailing zero), so force a text format before writing the
# value and then clear the format again so no stray number format is
# left behind on the cell.
$updates = @(
    @{ Addr = 'D2'; Value = '246.69' },
    @{ Addr = 'G2'; Value = '10' },
    @{ Addr = 'D3'; Value = '22.79' },
    @{ Addr = 'G3'; Value = '10' },
    @{ Addr = 'D4'; Value = '5.445' },
    @{ Addr = 'G4'; Value = '10' },
    @{ Addr = 'D5'; Value = '0.05769' },
    @{ Addr = 'G5'; Value = '10' },
    @{ Addr = 'D6'; Value = '3.429' },
    @{ Addr = 'G6'; Value = '10' },
    @{ Addr = 'B7'; Value = 'KuCoinToken' },
    @{ Addr = 'C7'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' },
    @{ Addr = 'D7'; Value = '6.323' },
    @{ Addr = 'E7'; Value = '6KuCoinTokenKCS' },
    @{ Addr = 'G7'; Value = '10' },
    @{ Addr = 'B8'; Value = 'MXToken' },
    @{ Addr = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Addr = 'D8'; Value = '0.8109' },
    @{ Addr = 'E8'; Value = '7MXTokenMX' },
    @{ Addr = 'G8'; Value = '10' },
    @{ Addr = 'B9'; Value = 'FTXToken' },
    @{ Addr = 'C9'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' },
    @{ Addr = 'D9'; Value = '0.8895' },
    @{ Addr = 'E9'; Value = '8FTXTokenFTT' },
    @{ Addr = 'G9'; Value = '10' },
    @{ Addr = 'B10'; Value = 'WazirX' },
    @{ Addr = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Addr = 'D10'; Value = '0.1447' },
    @{ Addr = 'E10'; Value = '9WazirXWRX' },
    @{ Addr = 'G10'; Value = '10' },
    @{ Addr = 'B11'; Value = 'MandalaExchangeToken' },
    @{ Addr = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Addr = 'D11'; Value = '0.07338' },
    @{ Addr = 'E11'; Value = '10MandalaExchangeTokenMDX' },
    @{ Addr = 'G11'; Value = '10' },
    @{ Addr = 'B12'; Value = 'BitrueCoin' },
    @{ Addr = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Addr = 'D12'; Value = '0.03134' },
    @{ Addr = 'E12'; Value = '11BitrueCoinBTRBestin24h' },
    @{ Addr = 'G12'; Value = '10' },
    @{ Addr = 'D13'; Value = '0.02966' },
    @{ Addr = 'G13'; Value = '10' },
    @{ Addr = 'G14'; Value = '10' },
    @{ Addr = 'D15'; Value = '3.928' },
    @{ Addr = 'G15'; Value = '10' },
    @{ Addr = 'D16'; Value = '0.001590' },
    @{ Addr = 'G16'; Value = '10' },
    @{ Addr = 'D17'; Value = '0.04798' },
    @{ Addr = 'G17'; Value = '10' },
    @{ Addr = 'D18'; Value = '0.0005847' },
    @{ Addr = 'G18'; Value = '10' },
    @{ Addr = 'D19'; Value = '0.006159' },
    @{ Addr = 'G19'; Value = '10' },
    @{ Addr = 'G20'; Value = '10' },
    @{ Addr = 'D21'; Value = '0.0009925' },
    @{ Addr = 'G21'; Value = '10' },
    @{ Addr = 'G22'; Value = '10' },
    @{ Addr = 'D23'; Value = '3.750' },
    @{ Addr = 'G23'; Value = '10' },
    @{ Addr = 'B24'; Value = 'BTSEToken' },
    @{ Addr = 'C24'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' },
    @{ Addr = 'D24'; Value = '2.199' },
    @{ Addr = 'E24'; Value = '23BTSETokenBTSE' },
    @{ Addr = 'G24'; Value = '10' },
    @{ Addr = 'B25'; Value = 'BitpandaEcosystemToken' },
    @{ Addr = 'C25'; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' },
    @{ Addr = 'D25'; Value = '0.3279' },
    @{ Addr = 'E25'; Value = '24BitpandaEcosystemTokenBEST' },
    @{ Addr = 'G25'; Value = '10' },
    @{ Addr = 'B26'; Value = 'ProBitToken' },
    @{ Addr = 'C26'; Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob' },
    @{ Addr = 'D26'; Value = '0.1319' },
    @{ Addr = 'E26'; Value = '25ProBitTokenPROB' },
    @{ Addr = 'G26'; Value = '10' },
    @{ Addr = 'D27'; Value = '0.0003157' },
    @{ Addr = 'G27'; Value = '10' },
    @{ Addr = 'G28'; Value = '10' },
    @{ Addr = 'G29'; Value = '10' },
    @{ Addr = 'G30'; Value = '10' },
    @{ Addr = 'G31'; Value = '10' },
    @{ Addr = 'G32'; Value = '10' },
    @{ Addr = 'G33'; Value = '10' },
    @{ Addr = 'G34'; Value = '10' },
    @{ Addr = 'G35'; Value = '10' },
    @{ Addr = 'G36'; Value = '10' },
    @{ Addr = 'G37'; Value = '10' },
    @{ Addr = 'G38'; Value = '10' },
    @{ Addr = 'G39'; Value = '10' },
    @{ Addr = 'D40'; Value = '0.03913' },
    @{ Addr = 'G40'; Value = '10' },
    @{ Addr = 'D41'; Value = '0.006765' },
    @{ Addr = 'G41'; Value = '10' },
    @{ Addr = 'D42'; Value = '0.1075' },
    @{ Addr = 'G42'; Value = '10' },
    @{ Addr = 'G43'; Value = '10' },
    @{ Addr = 'D44'; Value = '0.007095' },
    @{ Addr = 'G44'; Value = '10' },
    @{ Addr = 'D45'; Value = '0.00005641' },
    @{ Addr = 'G45'; Value = '10' },
    @{ Addr = 'G46'; Value = '10' },
    @{ Addr = 'D47'; Value = '0.3798' },
    @{ Addr = 'G47'; Value = '10' },
    @{ Addr = 'D48'; Value = '0.1682' },
    @{ Addr = 'G48'; Value = '10' },
    @{ Addr = 'D49'; Value = '0.00002099' },
    @{ Addr = 'G49'; Value = '10' },
    @{ Addr = 'G50'; Value = '10' },
    @{ Addr = 'G51'; Value = '10' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

$wb.Save()
